$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Size of crop image" column (J) values from 150 to 200
$ws.Range("J2").Value = 200
$ws.Range("J3").Value = 200
$ws.Range("J4").Value = 200

# Update the view: scroll so F1 is the top-left visible cell, and select H3
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H3").Select()
